$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Restore the original text of C4 (revert of the previous edit that replaced it
# with "Ghép nối các form").
$ws.Range("C4").Value = "Thiết kế giao diện cửa sổ chính (màn hình đăng nhập, các menu trỏ đến các module con,…). Xây dựng module quản lý người dùng"

# Update the active selection to match the reverted state.
$ws.Activate()
$ws.Range("E10").Select()
